$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the current row 47 (the old "2021-01-09" block).
# This shifts the old rows 47:50 down to 51:54 intact, and leaves fresh
# empty rows 47:50 to be filled with the new "2021-01-02" prediction block.
$ws.Rows("47:50").Insert()

# Row 47 - full row (Real/Prediction/difference + metrics), same as row 44's
# "Real" value carried forward, new Prediction/difference and recomputed metrics.
$ws.Range("A47").NumberFormat = "@"
$ws.Range("A47").Value = "2021-01-02"
$ws.Range("A47").ClearFormats()
$ws.Range("B47").Value = "03 Jan -- 09 Jan 2021"
$ws.Range("C47").Value = 94.56999999999999
$ws.Range("D47").Value = 343.5
$ws.Range("E47").Value = 248.93
$ws.Range("F47").Value = "KNN"
$ws.Range("G47").Value = 0.66
$ws.Range("H47").Value = 179.69
$ws.Range("I47").Value = 221.79
$ws.Range("J47").Value = 192.4
$ws.Range("K47").Value = 191.33

# Row 48
$ws.Range("A48").NumberFormat = "@"
$ws.Range("A48").Value = "2021-01-02"
$ws.Range("A48").ClearFormats()
$ws.Range("B48").Value = "10 Jan -- 16 Jan 2021"
$ws.Range("D48").Value = 331.74
$ws.Range("F48").Value = "KNN"

# Row 49
$ws.Range("A49").NumberFormat = "@"
$ws.Range("A49").Value = "2021-01-02"
$ws.Range("A49").ClearFormats()
$ws.Range("B49").Value = "17 Jan -- 23 Jan 2021"
$ws.Range("D49").Value = 321.16
$ws.Range("F49").Value = "KNN"

# Row 50
$ws.Range("A50").NumberFormat = "@"
$ws.Range("A50").Value = "2021-01-02"
$ws.Range("A50").ClearFormats()
$ws.Range("B50").Value = "24 Jan -- 30 Jan 2021"
$ws.Range("D50").Value = 319.23
$ws.Range("F50").Value = "KNN"
